$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 47576.5
$ws.Range("I21").Value = 40842.832
$ws.Range("J21").Value = 67777.5
$ws.Range("K21").Value = 40842.832
$ws.Range("L21").Value = 67777.5
$ws.Range("M21").Value = -40374.832
$ws.Range("N21").Value = -68713.5

$ws.Range("H23").Value = 47576.5
$ws.Range("I23").Value = 40842.832
$ws.Range("J23").Value = 67777.5
$ws.Range("K23").Value = 40842.832
$ws.Range("L23").Value = 67777.5
$ws.Range("M23").Value = -40608.832
$ws.Range("N23").Value = -68245.5

$ws.Range("H33").Value = 295.74194
$ws.Range("I33").Value = 188.51724
$ws.Range("J33").Value = 1850.5
$ws.Range("K33").Value = 188.51724
$ws.Range("L33").Value = 1850.5
$ws.Range("M33").Value = 40.48276000000001
$ws.Range("N33").Value = -2308.5

$ws.Range("H58").Value = 1636.6666
$ws.Range("J58").Value = 1920
$ws.Range("L58").Value = 5760
$ws.Range("N58").Value = -6060

$ws.Range("H80").Value = 20834742
$ws.Range("J80").Value = 1698.3
$ws.Range("L80").Value = 5094.9
$ws.Range("N80").Value = -7090.9

$ws.Range("H83").Value = 20834742
$ws.Range("J83").Value = 1698.3
$ws.Range("L83").Value = 15284.7
$ws.Range("N83").Value = -25268.7

$ws.Range("H94").Value = 976359.6
$ws.Range("I94").Value = 976359.6
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 976359.6
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -975908.6
$ws.Range("N94").ClearContents()

$ws.Range("H133").Value = 59918.8
$ws.Range("J133").Value = 59918.8
$ws.Range("L133").Value = 59918.8
$ws.Range("N133").Value = -70038.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()

$ws.Range("H102").Value = 27799.523
$ws.Range("I102").Value = 10475.883
$ws.Range("J102").Value = 101425
$ws.Range("K102").Value = 10475.883
$ws.Range("L102").Value = 101425
$ws.Range("M102").Value = -8853.883
$ws.Range("N102").Value = -104669

$ws.Range("H123").Value = 37714
$ws.Range("J123").Value = 37714
$ws.Range("L123").Value = 37714
$ws.Range("N123").Value = -47514

$ws.Range("H132").Value = 11112767
$ws.Range("I132").Value = 13514798
$ws.Range("K132").Value = 40544394
$ws.Range("M132").Value = -40541864

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

$ws.Range("H86").Value = 3796.8333
$ws.Range("I86").Value = 2570.5
$ws.Range("J86").Value = 6249.5
$ws.Range("K86").Value = 2570.5
$ws.Range("L86").Value = 6249.5
$ws.Range("M86").Value = -1447.5
$ws.Range("N86").Value = -8495.5

$ws.Range("H89").Value = 3796.8333
$ws.Range("I89").Value = 2570.5
$ws.Range("J89").Value = 6249.5
$ws.Range("K89").Value = 12852.5
$ws.Range("L89").Value = 31247.5
$ws.Range("M89").Value = -7236.5
$ws.Range("N89").Value = -42479.5

$ws.Range("H105").Value = 2324.88
$ws.Range("I105").Value = 2212.9583
$ws.Range("K105").Value = 2212.9583
$ws.Range("M105").Value = -465.9582999999998

$ws.Range("H125").Value = 50775
$ws.Range("J125").Value = 50775
$ws.Range("L125").Value = 50775
$ws.Range("N125").Value = -60615

$ws.Range("H134").Value = 2181.081
$ws.Range("I134").Value = 1654.7931
$ws.Range("K134").Value = 4964.379300000001
$ws.Range("M134").Value = -2429.379300000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 6174.08
$ws.Range("I3").Value = 2262.5
$ws.Range("J3").Value = 8014.8237
$ws.Range("K3").Value = 6787.5
$ws.Range("L3").Value = 24044.4711
$ws.Range("M3").Value = -6675.5
$ws.Range("N3").Value = -24268.4711

$ws.Range("H18").Value = 492.85715
$ws.Range("I18").Value = 525
$ws.Range("K18").Value = 1575
$ws.Range("M18").Value = -1406

$ws.Range("H62").Value = 3257
$ws.Range("I62").Value = 1500
$ws.Range("J62").Value = 5014
$ws.Range("K62").Value = 4500
$ws.Range("L62").Value = 15042
$ws.Range("M62").Value = -3814
$ws.Range("N62").Value = -16414

$ws.Range("H65").Value = 3257
$ws.Range("I65").Value = 1500
$ws.Range("J65").Value = 5014
$ws.Range("K65").Value = 13500
$ws.Range("L65").Value = 45126
$ws.Range("M65").Value = -10068
$ws.Range("N65").Value = -51990

$ws.Range("H68").Value = 1191.3625
$ws.Range("I68").Value = 802.88
$ws.Range("J68").Value = 1367.9454
$ws.Range("K68").Value = 2408.64
$ws.Range("L68").Value = 4103.8362
$ws.Range("M68").Value = -1597.64
$ws.Range("N68").Value = -5725.8362

$ws.Range("H71").Value = 1191.3625
$ws.Range("I71").Value = 802.88
$ws.Range("J71").Value = 1367.9454
$ws.Range("K71").Value = 7225.92
$ws.Range("L71").Value = 12311.5086
$ws.Range("M71").Value = -3169.92
$ws.Range("N71").Value = -20423.5086

$ws.Range("H80").Value = 38559940
$ws.Range("I80").Value = 2000000
$ws.Range("J80").Value = 40022340
$ws.Range("K80").Value = 6000000
$ws.Range("L80").Value = 120067020
$ws.Range("M80").Value = -5999064
$ws.Range("N80").Value = -120068892

$ws.Range("H83").Value = 38559940
$ws.Range("I83").Value = 2000000
$ws.Range("J83").Value = 40022340
$ws.Range("K83").Value = 18000000
$ws.Range("L83").Value = 360201060
$ws.Range("M83").Value = -17995320
$ws.Range("N83").Value = -360210420

$ws.Range("H113").Value = 2285.61
$ws.Range("I113").Value = 3058.8718
$ws.Range("J113").Value = 777.75
$ws.Range("K113").Value = 9176.615399999999
$ws.Range("L113").Value = 2333.25
$ws.Range("M113").Value = -7006.615399999999
$ws.Range("N113").Value = -6673.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 24995
$ws.Range("J26").Value = 24995
$ws.Range("L26").Value = 24995
$ws.Range("N26").Value = -25555

$ws.Range("H50").Value = 24995
$ws.Range("J50").Value = 24995
$ws.Range("L50").Value = 24995
$ws.Range("N50").Value = -25991

$ws.Range("H93").Value = 36986.832
$ws.Range("J93").Value = 36986.832
$ws.Range("L93").Value = 36986.832
$ws.Range("N93").Value = -40730.832

$ws.Range("H132").Value = 29414292
$ws.Range("I132").Value = 47620676
$ws.Range("K132").Value = 142862028
$ws.Range("M132").Value = -142859498

$ws.Range("H140").Value = 38033
$ws.Range("J140").Value = 38033
$ws.Range("L140").Value = 38033
$ws.Range("N140").Value = -48393

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H125").Value = 45715
$ws.Range("J125").Value = 45715
$ws.Range("L125").Value = 45715
$ws.Range("N125").Value = -55555

$ws.Range("H132").Value = 3669.3157
$ws.Range("I132").Value = 3213.3157
$ws.Range("J132").Value = 4125.316
$ws.Range("K132").Value = 9639.947100000001
$ws.Range("L132").Value = 12375.948
$ws.Range("M132").Value = -7109.947100000001
$ws.Range("N132").Value = -17435.948

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 40000
$ws.Range("J123").Value = 40000
$ws.Range("L123").Value = 40000
$ws.Range("N123").Value = -49800

$ws.Range("H137").Value = 60000
$ws.Range("J137").Value = 60000
$ws.Range("L137").Value = 60000
$ws.Range("N137").Value = -70200

$ws.Range("H138").Value = 44283.89
$ws.Range("J138").Value = 44283.89
$ws.Range("L138").Value = 44283.89
$ws.Range("N138").Value = -54563.89
